$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the two runs "SAT Oct 20" + " 11:58:44 IST 2018" into
# a single run by doing a Find/Replace over the exact same text; Word's
# replace collapses the matched range into one run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("SAT Oct 20 11:58:44 IST 2018", $false, $false, $false, $false, $false, $true, 1, $false, "SAT Oct 20 11:58:44 IST 2018", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: append a new "SAT Oct 21 ... MAMATHA CHICK IN" record right
# after the existing "SAT Oct 20" record's closing
# "Amount Received mode ... - CASH AND CLEARD" paragraph (that record is
# the one that starts right after the "SAT Oct 20" marker; the document
# has several other "- CASH AND CLEARD" lines earlier on, so we anchor
# the search to start just after the unique "SAT Oct 20" marker).
# ---------------------------------------------------------------------
$marker = $d.Content
$marker.Find.Execute("SAT Oct 20", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$tail = $d.Range($marker.End, $d.Content.End)
$tail.Find.Execute("- CASH AND CLEARD", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$insertAt = $tail.Paragraphs.First.Range
$insertAt.Collapse(0)

$lines = @(
  "",
  "SAT Oct 21 14:14:49 IST 2018",
  "Person Name`t`t`t`t- TRC",
  "Bill number`t`t`t`t- 8512",
  "---------------------------------------------------------------",
  "Item Name`t`t`t`t- CHOWCHOW",
  "Number of Pockets`t`t`t- 1",
  "Number of KGs`t`t`t- 67",
  "Rate`t`t`t`t`t- 10",
  "Total Price`t`t`t`t- 670.0",
  "Amount balance`t`t`t- 670.0",
  "",
  "Item Name`t`t`t`t- CARROT",
  "Number of Pockets`t`t`t- 1",
  "Number of KGs`t`t`t- 67",
  "Rate`t`t`t`t`t- 38",
  "Total Price`t`t`t`t- 2546.0",
  "Amount balance`t`t`t- 3216.0",
  "",
  "Item Name`t`t`t`t- CARROT",
  "Amount Received`t`t`t- 3216",
  "Amount Received mode`t`t- CASH AND CLEARD",
  ""
)

$insertAt.Text = ($lines -join "`r") + "`r"

# Colour the newly-inserted "Amount Received ... - 3216" line red, like
# the other "Amount Received" lines in the document.
$amt = $d.Content
$amt.Find.Execute("Amount Received`t`t`t- 3216", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$amt.Font.Color = 255
